$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://github.com/sangeetajoshi/E-Appointment_Final/tree/master/Source%20Code/EAppointment"

for ($r = 6; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
    $cell.HorizontalAlignment = -4131
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $false
}

$ws.Range("D13").Select() | Out-Null
